$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab "March" -> "März" (German)
$ws.Name = "März"

# Update the header row (row 4) to the German column headers
$ws.Range("A4").Value = "Email Addresse"
$ws.Range("B4").Value = "Nachname"
$ws.Range("C4").Value = "Vorname"

# Update the title cell (row 3) to the German translation
$ws.Range("A3").Value = "Webinar Anwesenheitsliste"

# Move the selection to A4 (matches the saved cursor position in the edit)
$ws.Range("A4").Select()
